$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 4 new rows for the new "2509" period block, pushing the
#    signature rows (old 24/25) down to 28/29.
$ws.Rows("20:23").Insert()

# 2. Update VALOR MORA total and Cant. Periodos count.
$ws.Range("E11").Value = 455520
$ws.Range("F13").Value = 2

# 3. Copy the formatting of the existing 4-worker block (rows 16-19,
#    which currently ends with the "closing" bottom border on row 19)
#    onto the freshly inserted rows 20-23, so row 23 becomes the new
#    closing row of the table.
$ws.Range("B16:J19").Copy()
$ws.Range("B20:J23").PasteSpecial(-4122)

# 4. Row 19 is no longer the last row of the table, so restyle it like
#    an interior row (matching rows 16-18).
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 5. Fill in the new rows: same four workers repeated for period 2509.
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1051446581"
$ws.Range("D20").Value = "ROQUE ALEXANDER PAJARO ACEVEDO"
$ws.Range("E20").Value = "2509"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1052080658"
$ws.Range("D21").Value = "GUILLERMO JOSE NARVAEZ ARIAS"
$ws.Range("E21").Value = "2509"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1006579007"
$ws.Range("D22").Value = "YEFERSON DAVID PEREZ CONTRERAS"
$ws.Range("E22").Value = "2509"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1002413113"
$ws.Range("D23").Value = "JHONATAN MEZA POLO"
$ws.Range("E23").Value = "2509"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

# 6. Center the "Periodo Mora" column for the whole worker table.
$ws.Range("E16:E23").HorizontalAlignment = -4108
